$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh timestamp applied to column D for every populated data row (2-39)
$newDate = 45997.373206018521
$ws.Range("D2:D39").Value = $newDate

# The report body (rows 19-56) is regenerated: old rows are cleared first
# (this also drops any shared strings that become unused), then the
# up-to-date rows 19-39 are written back in; rows 40-56 stay blank.
$ws.Range("A19:C56").ClearContents() | Out-Null

$data = @(
  @(19, "长沙特来电飞狐四方坪南区充电站", "406号直流", 45993.542002314818),
  @(20, "长沙特来电飞狐四方坪南区充电站", "201号直流", 45994.55159722222),
  @(21, "长沙特来电飞狐四方坪东区充电站", "005A号直流", 45995.092395833337),
  @(22, "长沙特来电飞狐四方坪东区充电站", "003B号直流", 45995.604780092595),
  @(23, "长沙特来电飞狐四方坪西区充电站", "705号直流", 45995.667939814812),
  @(24, "长沙特来电飞狐四方坪东区充电站", "001B号直流", 45996.17591435185),
  @(25, "长沙特来电飞狐四方坪西区充电站", "903号直流", 45996.410590277781),
  @(26, "长沙特来电飞狐四方坪南区充电站", "104号直流", 45996.450902777775),
  @(27, "长沙特来电飞狐四方坪西区充电站", "405号直流", 45996.549861111111),
  @(28, "长沙特来电飞狐四方坪西区充电站", "503号直流", 45996.554583333331),
  @(29, "长沙市开福区高岭香江国际城充电站建设项目", "108号直流", 45996.562662037039),
  @(30, "长沙特来电飞狐四方坪东区充电站", "102号直流", 45996.572453703702),
  @(31, "长沙特来电飞狐四方坪西区充电站", "B01号直流", 45996.593472222223),
  @(32, "长沙特来电飞狐四方坪西区充电站", "A01号直流", 45996.609143518515),
  @(33, "长沙特来电飞狐四方坪东区充电站", "002B号直流", 45996.618958333333),
  @(34, "长沙特来电飞狐四方坪南区充电站", "305号直流", 45996.621736111112),
  @(35, "长沙市开福区高岭香江国际城充电站建设项目", "110号直流", 45996.642685185187),
  @(36, "长沙特来电飞狐四方坪南区充电站", "203号直流", 45996.685231481482),
  @(37, "长沙市开福区高岭香江国际城充电站建设项目", "202号直流", 45996.794722222221),
  @(38, "长沙市开福区高岭香江国际城充电站建设项目", "203号直流", 45996.817743055559),
  @(39, "长沙特来电飞狐四方坪西区充电站", "A04号直流", 45996.835162037038)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Range("A$r").Value = $row[1]
  $ws.Range("B$r").Value = $row[2]
  $ws.Range("C$r").Value = $row[3]
}

# Selection moved to E17 in the saved view
$ws.Range("E17").Select() | Out-Null
